# MSME Country Indicators - Iran, Islamic Rep. - Summary sheet
# Refresh a handful of percentage figures with more precise decimal values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Helper: a cell whose current format/style we can borrow so that writing a
# numeric-looking string back doesn't leave the cell re-styled (Excel
# auto-detects "19.78" etc. as a number unless the cell is Text-formatted,
# which otherwise bumps the cell's style index).
$templateStyle = $ws.Range("B12").Style

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $templateStyle
}

# Enterprises density (per 1000 people)
Set-TextValue "B13" "19.78"
Set-TextValue "C13" "0.29"
Set-TextValue "D13" "20.07"

# Employment (% of total) -- SMEs column (C14, 42.9) is unchanged
Set-TextValue "B14" "11.74"
Set-TextValue "D14" "54.64"

# Enterprises (% of total)
Set-TextValue "B16" "98.43"
Set-TextValue "C16" "1.47"
Set-TextValue "D16" "99.89"
